$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verifications")

# Append the two new verification upload items for "847 scenario 1"
$ws.Range("A96").Value = "Upload a copy of the LPA planning decision showing approval of the outline planning permission"
$ws.Range("A97").Value = "Upload a copy of the LPA decision notice that you are appealing against"
